# Update epexspot_prices.xlsx with the latest daily prices.
# 1) "Prix Spot" sheet: append a new date column (AG) with hourly prices for 16-jul.
# 2) "Gaz" sheet: append a new row for 2025-07-14.
# 3) "CO2" sheet: append a new row for 2025-07-14.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" -- add column AG (16-jul)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (AF1) onto the new header
# cell (AG1) so the new column matches the existing bold/bordered/centered
# header style, then set its text.
$ws1.Range("AF1").Copy()
$ws1.Cells.Item(1, 33).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Cells.Item(1, 33).Value = "16-jul"

$pricesAG = @(
    74.97,
    56.53,
    55.26,
    48.19,
    44.71,
    46.97,
    52.33,
    62.55,
    75.2,
    66.98999999999999,
    42.63,
    28.81,
    53.6,
    36.26,
    33.25,
    29.01,
    39.76,
    57.57,
    79.09999999999999,
    109.62,
    122.84,
    108.4,
    117.4,
    111.14
)

for ($i = 0; $i -lt $pricesAG.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 33).Value = $pricesAG[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" -- add row 30 (2025-07-14)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

$dateCell2 = $ws2.Cells.Item(30, 1)
# Force text storage so "2025-07-14" isn't auto-converted to a date serial.
$dateCell2.NumberFormat = "@"
$dateCell2.Value = "2025-07-14"
# Reset formatting to match the plain (unstyled) cells above it.
$ws2.Range("A29").Copy()
$dateCell2.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Cells.Item(30, 2).Value = 34.275

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" -- add row 30 (2025-07-14)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$dateCell3 = $ws3.Cells.Item(30, 1)
$dateCell3.NumberFormat = "@"
$dateCell3.Value = "2025-07-14"
$ws3.Range("A29").Copy()
$dateCell3.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Cells.Item(30, 2).Value = 69.59999999999999
